# ---------------------------------------------------------------------------
# Applies the two changes captured by the commit:
#   1. The table on slide 5 switches from the deck's custom "Table_0" style
#      ({35C023E5-1390-4394-8B6E-CA9D61337C18}, defined in tableStyles.xml)
#      to PowerPoint's built-in "No Style, Table Grid" style
#      ({7C8D13BF-48A7-4531-BA89-1CE60BAFDF38}).
#   2. The deck's theme palette ("Integral" / "Red Violet") is replaced by
#      the default Office theme palette ("Office").
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 5 (graphicFrame is Shapes.Item(2)) -----
$slide5 = $p.Slides.Item(5)
$tableShape = $null
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $candidate = $slide5.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
        break
    }
}
$table = $tableShape.Table
$table.ApplyStyle("{7C8D13BF-48A7-4531-BA89-1CE60BAFDF38}", $false)

# --- 2. Swap the "Integral"/Red Violet theme colours for the Office theme --
# Theme colour slots, in MsoThemeColorSchemeIndex order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeColors = @(
    0,            # dk1      000000
    16777215,     # lt1      FFFFFF
    6968388,      # dk2      44546A
    15132391,     # lt2      E7E6E6
    13998939,     # accent1  5B9BD5
    3243501,      # accent2  ED7D31
    10855845,     # accent3  A5A5A5
    49407,        # accent4  FFC000
    12874308,     # accent5  4472C4
    4697456,      # accent6  70AD47
    12673797,     # hlink    0563C1
    7491477       # folHlink 954F72
)

$colorScheme = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Colors($i).RGB = $officeColors[$i - 1]
}
